$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82, shifting the existing rows (82..187) down to (83..188).
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row with the new weekly price-report record.
$ws.Cells.Item(82, 1).Value = 4
$ws.Cells.Item(82, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(82, 3).Value = "Los Lagos"
$ws.Cells.Item(82, 4).Value = 44539
$ws.Cells.Item(82, 5).Value = 10
$ws.Cells.Item(82, 6).Value = 100112003
$ws.Cells.Item(82, 7).Value = "Ajo"
$ws.Cells.Item(82, 8).Value = "Chino"
$ws.Cells.Item(82, 9).Value = "Primera"
$ws.Cells.Item(82, 10).Value = 60
$ws.Cells.Item(82, 11).Value = 21000
$ws.Cells.Item(82, 12).Value = 22000
$ws.Cells.Item(82, 13).Value = 21500
$ws.Cells.Item(82, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(82, 15).Value = "China"
$ws.Cells.Item(82, 16).Value = 2150
$ws.Cells.Item(82, 17).Value = 10
$ws.Cells.Item(82, 18).Value = "Hortaliza"
